# Corrigir Regras de Verificação
#
# 1) "RF define ..."  -> "Requisitos Funcionais (RF) define ..."
# 2) "RNF irá ..."     -> "Requisitos Não Funcionais (RNF) irá ..."
#    (the document's "_GoBack" bookmark is relocated from the trailing
#    empty paragraph to sit between "(" and "RNF)" in this sentence)

$d = $word.ActiveDocument

# --- 1) Expand "RF" into "Requisitos Funcionais (RF)" ----------------------
$rngRF = $d.Content
$rngRF.Find.Execute("RF ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rfRange = $d.Range($rngRF.Start, $rngRF.Start + 2)
$rfRange.Text = "Requisitos Funcionais (RF)"

# --- 2) Expand "RNF" into "Requisitos Não Funcionais (RNF)" ----------------
$rngRNF = $d.Content
$rngRNF.Find.Execute("RNF ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rnfStart = $rngRNF.Start
$rnfRange = $d.Range($rnfStart, $rnfStart + 3)
$rnfRange.Text = "Requisitos Não Funcionais (RNF)"

# Move the "_GoBack" bookmark so it sits right after the opening
# parenthesis, i.e. between "(" and "RNF)". Re-adding a bookmark with the
# same name relocates it, automatically removing it from its old spot
# (the trailing empty paragraph), which is left empty afterwards.
$prefixLen = "Requisitos Não Funcionais (".Length
$bmPos = $rnfStart + $prefixLen
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
